{"js": "// Bump the version string \"V3.3.6.2\" -> \"V3.3.7.2\" in the title block,\n// reproducing the author's edit exactly: the four runs that make up the\n// version text (\"V3.3.\" / \"6\"->\"7\" / \".\" / \"2\") are rebuilt *after* the\n// existing \"_GoBack\" bookmark (previously they sat before it).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that holds the version string (e.g. \"V3.3.6.2\").\nparas.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst versionRe = /^V(\\d+)\\.(\\d+)\\.(\\d+)\\.(\\d+)$/;\nlet target = null;\nlet m = null;\nfor (const p of paras.items) {\n  const t = (p.text || \"\").trim();\n  const mm = t.match(versionRe);\n  if (mm) {\n    target = p;\n    m = mm;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'V#.#.#.#' version paragraph\");\n}\n\n// Increment the third version component (3.3.6.2 -> 3.3.7.2), matching\n// the commit's 6 -> 7 bump.\nconst major = m[1];\nconst minor = m[2];\nconst patch = String(Number(m[3]) + 1);\nconst build = m[4];\n\nconst flatOpc = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:pPr><w:autoSpaceDE w:val=\"0\"/><w:autoSpaceDN w:val=\"0\"/><w:adjustRightInd w:val=\"0\"/><w:ind w:right=\"340\"/><w:jc w:val=\"center\"/><w:rPr><w:rFonts w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:b/><w:bCs/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr></w:pPr>\n<w:bookmarkStart w:id=\"4\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"4\"/>\n<w:r><w:rPr><w:rFonts w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>${\"V\" + major + \".\" + minor + \".\"}</w:t></w:r>\n<w:r><w:rPr><w:rFonts w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>${patch}</w:t></w:r>\n<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\" w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>.</w:t></w:r>\n<w:r><w:rPr><w:rFonts w:hint=\"default\" w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>${build}</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Bump the version string \"V3.3.6.2\" -> \"V3.3.7.2\" in the title block,\n# reproducing the author's edit exactly: the four runs that make up the\n# version text (\"V3.3.\" / \"6\"->\"7\" / \".\" / \"2\") are rebuilt *after* the\n# existing \"_GoBack\" bookmark (previously they sat before it).\n$d = $word.ActiveDocument\n\n$target = $null\n$m = $null\nforeach ($para in $d.Paragraphs) {\n  $t = $para.Range.Text\n  $trimmed = $t.TrimEnd([char]13, [char]7)\n  if ($trimmed -match '^V(\\d+)\\.(\\d+)\\.(\\d+)\\.(\\d+)$') {\n    $target = $para\n    $m = $Matches\n    break\n  }\n}\n\nif ($target -eq $null) {\n  throw \"Could not find the 'V#.#.#.#' version paragraph\"\n}\n\n$majorMinor = \"V\" + $m[1] + \".\" + $m[2] + \".\"\n$patch = [string]([int]$m[3] + 1)\n$build = $m[4]\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:autoSpaceDE w:val=\"0\"/><w:autoSpaceDN w:val=\"0\"/><w:adjustRightInd w:val=\"0\"/><w:ind w:right=\"340\"/><w:jc w:val=\"center\"/><w:rPr><w:rFonts w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:b/><w:bCs/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr></w:pPr><w:bookmarkStart w:id=\"4\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"4\"/><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>' + $majorMinor + '</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>' + $patch + '</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\" w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=\"default\" w:ascii=\"Helvetica\" w:hAnsi=\"Helvetica\" w:eastAsia=\"PingFang SC\" w:cs=\"Helvetica\"/><w:kern w:val=\"0\"/><w:szCs w:val=\"21\"/></w:rPr><w:t>' + $build + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$target.Range.InsertXML($xml)\n"}
